# Update "想去人数" (want-to-go count) figures in column F across the four
# sheets to reflect the newly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value  = 283
$ws.Range("F19").Value = 1653
$ws.Range("F27").Value = 4354
$ws.Range("F28").Value = 298

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value  = 608
$ws.Range("F27").Value = 6230

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value  = 1956
$ws.Range("F10").Value = 1269
$ws.Range("F13").Value = 2070
$ws.Range("F14").Value = 8774

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 1956
$ws.Range("F9").Value  = 1269
$ws.Range("F16").Value = 608
$ws.Range("F31").Value = 298
$ws.Range("F49").Value = 6230

$wb.Save()
